$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set each changed cell as literal text (preserving formatting/leading
# zeros/percent signs exactly) by forcing Text number format, assigning the
# value, then restoring the default "Normal" style so no stray per-cell
# style index is left behind.
$cellValues = @{
    'D2' = '315.94'
    'D3' = '39.41'
    'E3' = '-0.62%'
    'D4' = '5.127'
    'E4' = '-0.36%'
    'D5' = '0.08176'
    'E5' = '1.01%'
    'D6' = '1.981'
    'E6' = '1.90%'
    'B7' = 'KuCoinToken'
    'C7' = 'https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs'
    'D7' = '8.319'
    'E7' = '2.09%'
    'B8' = 'MXToken'
    'C8' = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
    'D8' = '0.9381'
    'E8' = '1.13%'
    'B9' = 'LiechtensteinCryptoassetsExchange'
    'C9' = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
    'D9' = '0.1301'
    'E9' = '-7.08%'
    'B10' = 'WazirX'
    'C10' = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
    'D10' = '0.1972'
    'E10' = '2.95%'
    'B11' = 'MandalaExchangeToken'
    'C11' = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
    'D11' = '0.09046'
    'E11' = '-0.53%'
    'B12' = 'BitrueCoin'
    'C12' = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
    'D12' = '0.03498'
    'E12' = '0.17%'
    'B13' = 'BitMartToken'
    'C13' = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
    'D13' = '0.09747'
    'E13' = '-0.73%'
    'B14' = 'BitForexToken'
    'C14' = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
    'D14' = '0.001414'
    'E14' = '1.29%'
    'B15' = 'TigerCash'
    'C15' = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
    'D15' = '0.006014'
    'E15' = '0.63%'
    'B16' = 'LEO'
    'C16' = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
    'D16' = '3.641'
    'E16' = '-7.69%'
    'B17' = 'GateToken'
    'C17' = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
    'D17' = '4.369'
    'E17' = '3.27%'
    'E18' = '-1.31%'
    'E19' = '1.86%'
    'D20' = '0.1301'
    'E20' = '-3.45%'
    'D21' = '4.956'
    'E21' = '6.43%'
    'D22' = '0.2581'
    'E22' = '6.55%'
    'D23' = '0.04362'
    'E23' = '-0.29%'
    'D24' = '0.001243'
    'E24' = '1.36%'
    'D25' = '0.004764'
    'E25' = '9.72%'
    'D26' = '0.0003892'
    'E26' = '199.07%'
    'E27' = '-7.58%'
    'D39' = '0.02215'
    'E39' = '9.02%'
    'E40' = '2.75%'
    'D41' = '0.007755'
    'E41' = '4.73%'
    'D42' = '0.01027'
    'E42' = '6.03%'
    'D43' = '0.1398'
    'E43' = '2.69%'
    'D44' = '0.002101'
    'E44' = '-1.46%'
    'D45' = '0.009273'
    'E45' = '6.57%'
    'D46' = '0.00006939'
    'E46' = '9.05%'
    'D47' = '0.00000000750'
    'E47' = '-0.04%'
    'D48' = '0.002884'
    'E48' = '0.64%'
    'D49' = '0.001691'
    'E49' = '30.02%'
    'D50' = '0.00002101'
    'E50' = '-0.04%'
    'D51' = '0.0002001'
    'E51' = '-0.04%'
}

foreach ($addr in $cellValues.Keys) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $cellValues[$addr]
    $rng.Style = "Normal"
}
